$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 228
$ws1.Range("F5").Value = 9072
$ws1.Range("F6").Value = 531
$ws1.Range("F8").Value = 148
$ws1.Range("F9").Value = 204
$ws1.Range("F10").Value = 327
$ws1.Range("F11").Value = 376
$ws1.Range("F14").Value = 409
$ws1.Range("F15").Value = 11604
$ws1.Range("F16").Value = 11604
$ws1.Range("F26").Value = 213
$ws1.Range("F27").Value = 34
$ws1.Range("F28").Value = 17
$ws1.Range("F29").Value = 144
$ws1.Range("F30").Value = 2702
$ws1.Range("F33").Value = 2089
$ws1.Range("F34").Value = 54
$ws1.Range("F35").Value = 49
$ws1.Range("F36").Value = 2125
$ws1.Range("F37").Value = 954
$ws1.Range("F38").Value = 4150
$ws1.Range("F39").Value = 313
$ws1.Range("F40").Value = 3044
$ws1.Range("F41").Value = 1282
$ws1.Range("F43").Value = 82
$ws1.Range("F44").Value = 381
$ws1.Range("F45").Value = 424
$ws1.Range("F48").Value = 173
$ws1.Range("F50").Value = 111

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 11
$ws2.Range("F9").Value = 38
$ws2.Range("F20").Value = 69

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 11
$ws4.Range("F7").Value = 228
$ws4.Range("F9").Value = 9072
$ws4.Range("F10").Value = 531
$ws4.Range("F11").Value = 38
$ws4.Range("F13").Value = 204
$ws4.Range("F14").Value = 327
$ws4.Range("F15").Value = 376
$ws4.Range("F18").Value = 11604
$ws4.Range("F28").Value = 213
$ws4.Range("F29").Value = 34
$ws4.Range("F30").Value = 17
$ws4.Range("F32").Value = 144
$ws4.Range("F34").Value = 2089
$ws4.Range("F35").Value = 54
$ws4.Range("F36").Value = 49
$ws4.Range("F37").Value = 2125
$ws4.Range("F38").Value = 954
$ws4.Range("F40").Value = 4150
$ws4.Range("F41").Value = 313
$ws4.Range("F42").Value = 3044
$ws4.Range("F43").Value = 1282
$ws4.Range("F44").Value = 82
$ws4.Range("F45").Value = 381
$ws4.Range("F48").Value = 173
$ws4.Range("F50").Value = 111
